$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "3 hours– creating a unity sample project to refresh Unity skill, by watch a simple unity tutorial  "
$ws.Range("B13").Value = "6 hours – trying to create a first person charcter movement in 3d and watch tutorial video online"
$ws.Range("B14").Value = "3 hours – have a group meeting with our supervisor to talk about the specific detail of the project.`n 3 hour – playing around with unity probuilder."
$ws.Range("B16").Value = "2 hours writing the design document."
$ws.Range("B19").Value = "2 hour - have a group meeting to talk about the project"
$ws.Range("B20").Value = "4 hour - watch tutorial on creating multiplayer gaming using FishNet `n2 hour – try using Fishnet to run a multiplayer for the sample game. The result end up being that you can only run multiplayer on the same computer making it not truly multiplayer`n2 hour – create a sample game for testing multiplayer "
$ws.Range("B21").Value = "1 hour – watching somemore random youtube video on Unity.`n1 hour - watching a tutorial on Photon Fusion. `n"
$ws.Range("B22").Value = "6 hour-find documentation and  watching youtube video on creating a first person controller in unity "
$ws.Range("B24").Value = "1 hour 15 minutes - looking for and adding player character asset and testing"
$ws.Range("B26").Value = "7 hour 30 minutes - creating the first playable first personcharacter by adding in basic movement, camera view, animation, and gunscript for basic shooting"
$ws.Range("B27").Value = "1 hour - rewatch youtube video on how to use Photon Pun to set up a multiplayer network`n1 hour - testing multiplayer with alda`n5 hour - trying to implement multiplayer feature to the game testing and fixing bug`n  "
$ws.Range("B28").Value = "3 hour - watching video on animation rigging in unity "
$ws.Range("B29").Value = "1 hour- group meeting about our current progress"
$ws.Range("B30").Value = "2 hours- testing new design ideas for player character"
$ws.Range("B33").Value = "5 hour - creating a new character design from scratch"
$ws.Range("B35").Value = "5 hour -creating a new weapon script and weapon system and adding side project to the repository`n1 hour - have a group meeting on our progress"
$ws.Range("B36").Value = "7 hour - creating a weapon switching system and script for multiple new weapon. Testing and creating and adding custom idle, fire, reload and loadout animation for each of them. Also getting them ready for multiplayer"
$ws.Range("B37").Value = "3 hour - fixing bug with the gaming and add health and ammo bar for the player character.`n2 hour - learning about RPC for Photon Pun and adding it to the game"
$ws.Range("B40").Value = "2 hour -trying to recreate the time log and fixing github problem that is not letting update timelog"

# B19 uses a non-wrapping center alignment style (new cellXfs entry)
$ws.Range("B19").WrapText = $false

# Update selection to B40 to match final state
$ws.Range("B40").Select()
